$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Cells.Item(38, 8).Value = 223.66667
$ws.Cells.Item(38, 9).Value = 68.40000000000001
$ws.Cells.Item(38, 10).Value = 1000
$ws.Cells.Item(38, 11).Value = 205.2
$ws.Cells.Item(38, 12).Value = 3000
$ws.Cells.Item(38, 13).Value = 166.8
$ws.Cells.Item(38, 14).Value = -3744

# Row 53
$ws.Cells.Item(53, 8).Value = 187.24243
$ws.Cells.Item(53, 9).Value = 135.65384
$ws.Cells.Item(53, 10).Value = 378.85715
$ws.Cells.Item(53, 11).Value = 135.65384
$ws.Cells.Item(53, 12).Value = 378.85715
$ws.Cells.Item(53, 13).Value = 501.34616
$ws.Cells.Item(53, 14).Value = -1652.85715

# Row 98
$ws.Cells.Item(98, 8).Value = 2572.1428
$ws.Cells.Item(98, 9).Value = 2500.7693
$ws.Cells.Item(98, 10).Value = 3500
$ws.Cells.Item(98, 11).Value = 2500.7693
$ws.Cells.Item(98, 12).Value = 3500
$ws.Cells.Item(98, 13).Value = -1002.7693
$ws.Cells.Item(98, 14).Value = -6496

# Row 118
$ws.Cells.Item(118, 8).Value = 751.46155
$ws.Cells.Item(118, 9).Value = 522.4167
$ws.Cells.Item(118, 10).Value = 3500
$ws.Cells.Item(118, 11).Value = 1567.2501
$ws.Cells.Item(118, 12).Value = 10500
$ws.Cells.Item(118, 13).Value = 89.74990000000003
$ws.Cells.Item(118, 14).Value = -13814

# Row 122
$ws.Cells.Item(122, 8).Value = 2572.1428
$ws.Cells.Item(122, 9).Value = 2500.7693
$ws.Cells.Item(122, 10).Value = 3500
$ws.Cells.Item(122, 11).Value = 7502.3079
$ws.Cells.Item(122, 12).Value = 10500
$ws.Cells.Item(122, 13).Value = -5052.3079
$ws.Cells.Item(122, 14).Value = -15400

# Row 123
$ws.Cells.Item(123, 8).Value = 35000
$ws.Cells.Item(123, 10).Value = 35000
$ws.Cells.Item(123, 12).Value = 35000
$ws.Cells.Item(123, 14).Value = -44800

# Row 124
$ws.Cells.Item(124, 8).Value = 35250
$ws.Cells.Item(124, 10).Value = 35250
$ws.Cells.Item(124, 12).Value = 35250
$ws.Cells.Item(124, 14).Value = -45070

# Row 125
$ws.Cells.Item(125, 8).Value = 1400
$ws.Cells.Item(125, 9).Value = 2000
$ws.Cells.Item(125, 10).Value = 1314.2858
$ws.Cells.Item(125, 11).Value = 18000
$ws.Cells.Item(125, 12).Value = 11828.5722
$ws.Cells.Item(125, 13).Value = -15540
$ws.Cells.Item(125, 14).Value = -16748.5722

# Row 127
$ws.Cells.Item(127, 8).Value = 2620.0227
$ws.Cells.Item(127, 9).Value = 1321.2222
$ws.Cells.Item(127, 10).Value = 2954
$ws.Cells.Item(127, 11).Value = 3963.6666
$ws.Cells.Item(127, 12).Value = 8862
$ws.Cells.Item(127, 13).Value = 996.3334000000004
$ws.Cells.Item(127, 14).Value = -18782

# Row 129
$ws.Cells.Item(129, 8).Value = 1029960
$ws.Cells.Item(129, 9).Value = 464.8
$ws.Cells.Item(129, 10).Value = 1196007.6
$ws.Cells.Item(129, 11).Value = 1394.4
$ws.Cells.Item(129, 12).Value = 3588022.8
$ws.Cells.Item(129, 13).Value = 3605.6
$ws.Cells.Item(129, 14).Value = -3598022.8

# Row 137
$ws.Cells.Item(137, 8).Value = 1157.7715
$ws.Cells.Item(137, 9).Value = 977.3333
$ws.Cells.Item(137, 10).Value = 2240.4
$ws.Cells.Item(137, 11).Value = 2931.9999
$ws.Cells.Item(137, 12).Value = 6721.200000000001
$ws.Cells.Item(137, 13).Value = -381.9998999999998
$ws.Cells.Item(137, 14).Value = -11821.2

# Row 138
$ws.Cells.Item(138, 8).Value = 4208.4
$ws.Cells.Item(138, 9).Value = 2469.3684
$ws.Cells.Item(138, 10).Value = 5126.222
$ws.Cells.Item(138, 11).Value = 7408.1052
$ws.Cells.Item(138, 12).Value = 15378.666
$ws.Cells.Item(138, 13).Value = -2268.1052
$ws.Cells.Item(138, 14).Value = -25658.666


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Cells.Item(5, 8).Value = 86.333336
$ws.Cells.Item(5, 9).Value = 86.333336
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 86.333336
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).ClearContents()
$ws.Cells.Item(5, 14).Value = 25.666664

# Row 32
$ws.Cells.Item(32, 8).Value = 33531.594
$ws.Cells.Item(32, 9).Value = 25205.885
$ws.Cells.Item(32, 11).Value = 25205.885
$ws.Cells.Item(32, 13).Value = -24918.885

# Row 45
$ws.Cells.Item(45, 8).Value = 1611.85
$ws.Cells.Item(45, 9).Value = 1811.1818
$ws.Cells.Item(45, 10).Value = 1368.2222
$ws.Cells.Item(45, 11).Value = 1811.1818
$ws.Cells.Item(45, 12).Value = 1368.2222
$ws.Cells.Item(45, 13).Value = -1434.1818
$ws.Cells.Item(45, 14).Value = -2122.2222

# Row 74
$ws.Cells.Item(74, 8).Value = 1132.3658
$ws.Cells.Item(74, 9).Value = 1062.8918
$ws.Cells.Item(74, 10).Value = 1775
$ws.Cells.Item(74, 11).Value = 1062.8918
$ws.Cells.Item(74, 12).Value = 1775
$ws.Cells.Item(74, 13).Value = -188.8918000000001
$ws.Cells.Item(74, 14).Value = -3523

# Row 77
$ws.Cells.Item(77, 8).Value = 1132.3658
$ws.Cells.Item(77, 9).Value = 1062.8918
$ws.Cells.Item(77, 10).Value = 1775
$ws.Cells.Item(77, 11).Value = 5314.459000000001
$ws.Cells.Item(77, 12).Value = 8875
$ws.Cells.Item(77, 13).Value = -946.4590000000007
$ws.Cells.Item(77, 14).Value = -17611

# Row 132
$ws.Cells.Item(132, 8).Value = 2279.2424
$ws.Cells.Item(132, 9).Value = 1738.4166
$ws.Cells.Item(132, 10).Value = 3721.4443
$ws.Cells.Item(132, 11).Value = 5215.2498
$ws.Cells.Item(132, 12).Value = 11164.3329
$ws.Cells.Item(132, 13).Value = -2685.2498
$ws.Cells.Item(132, 14).Value = -16224.3329


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Cells.Item(4, 8).Value = 86.333336
$ws.Cells.Item(4, 9).Value = 86.333336
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 86.333336
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).ClearContents()
$ws.Cells.Item(4, 14).Value = 28.666664

# Row 134
$ws.Cells.Item(134, 8).Value = 77578.37
$ws.Cells.Item(134, 9).Value = 3663.7273
$ws.Cells.Item(134, 10).Value = 402802.8
$ws.Cells.Item(134, 11).Value = 10991.1819
$ws.Cells.Item(134, 12).Value = 1208408.4
$ws.Cells.Item(134, 13).Value = -8456.1819
$ws.Cells.Item(134, 14).Value = -1213478.4


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 2121.5334
$ws.Cells.Item(16, 9).Value = 1202.2
$ws.Cells.Item(16, 10).Value = 2581.2
$ws.Cells.Item(16, 11).Value = 1202.2
$ws.Cells.Item(16, 12).Value = 2581.2
$ws.Cells.Item(16, 13).Value = -915.2
$ws.Cells.Item(16, 14).Value = -3155.2

# Row 113
$ws.Cells.Item(113, 8).Value = 2121.5334
$ws.Cells.Item(113, 9).Value = 1202.2
$ws.Cells.Item(113, 10).Value = 2581.2
$ws.Cells.Item(113, 11).Value = 1202.2
$ws.Cells.Item(113, 12).Value = 2581.2
$ws.Cells.Item(113, 13).Value = 967.8
$ws.Cells.Item(113, 14).Value = -6921.2

# Row 132
$ws.Cells.Item(132, 8).Value = 1456.862
$ws.Cells.Item(132, 9).Value = 1065.2727
$ws.Cells.Item(132, 10).Value = 2687.5715
$ws.Cells.Item(132, 11).Value = 3195.8181
$ws.Cells.Item(132, 12).Value = 8062.7145
$ws.Cells.Item(132, 13).Value = -665.8181
$ws.Cells.Item(132, 14).Value = -13122.7145


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Cells.Item(5, 8).Value = 1935.7307
$ws.Cells.Item(5, 9).Value = 1291.7646
$ws.Cells.Item(5, 10).Value = 2248.5144
$ws.Cells.Item(5, 11).Value = 3875.2938
$ws.Cells.Item(5, 12).Value = 6745.5432
$ws.Cells.Item(5, 13).Value = -3763.2938
$ws.Cells.Item(5, 14).Value = -6969.5432

# Row 80
$ws.Cells.Item(80, 8).Value = 6941.75
$ws.Cells.Item(80, 10).Value = 6877.778
$ws.Cells.Item(80, 12).Value = 20633.334
$ws.Cells.Item(80, 14).Value = -22505.334

# Row 83
$ws.Cells.Item(83, 8).Value = 6941.75
$ws.Cells.Item(83, 10).Value = 6877.778
$ws.Cells.Item(83, 12).Value = 61900.002
$ws.Cells.Item(83, 14).Value = -71260.00200000001

# Row 113
$ws.Cells.Item(113, 8).Value = 832.8125
$ws.Cells.Item(113, 9).Value = 994.4286
$ws.Cells.Item(113, 10).Value = 707.1111
$ws.Cells.Item(113, 11).Value = 2983.2858
$ws.Cells.Item(113, 12).Value = 2121.3333
$ws.Cells.Item(113, 13).Value = -813.2857999999997
$ws.Cells.Item(113, 14).Value = -6461.3333

# Row 132
$ws.Cells.Item(132, 8).Value = 1372.5
$ws.Cells.Item(132, 9).Value = 1401
$ws.Cells.Item(132, 10).Value = 1361.5385
$ws.Cells.Item(132, 11).Value = 12609
$ws.Cells.Item(132, 12).Value = 12253.8465
$ws.Cells.Item(132, 13).Value = -10079
$ws.Cells.Item(132, 14).Value = -17313.8465

# Row 135
$ws.Cells.Item(135, 8).Value = 1935.7307
$ws.Cells.Item(135, 9).Value = 1291.7646
$ws.Cells.Item(135, 10).Value = 2248.5144
$ws.Cells.Item(135, 11).Value = 11625.8814
$ws.Cells.Item(135, 12).Value = 20236.6296
$ws.Cells.Item(135, 13).Value = -9090.8814
$ws.Cells.Item(135, 14).Value = -25306.6296


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Cells.Item(132, 8).Value = 3477.9736
$ws.Cells.Item(132, 9).Value = 2329.625
$ws.Cells.Item(132, 10).Value = 5446.5713
$ws.Cells.Item(132, 11).Value = 6988.875
$ws.Cells.Item(132, 12).Value = 16339.7139
$ws.Cells.Item(132, 13).Value = -4458.875
$ws.Cells.Item(132, 14).Value = -21399.7139

